$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric identifiers
$ws.Range("A2").Value = 80139019
$ws.Range("B2").Value = 88921
$ws.Range("E2").Value = 5741

# Species names
$ws.Range("F2").Value = "Tjockfotad fingersvamp"
$ws.Range("G2").Value = "Ramaria flavescens"
$ws.Range("H2").Value = "(Schaeff.) R. H. Petersen"

# Antal (I2) cleared, Enhet (J2), Alder-Stadium (K2), Kon (L2), Metod (N2) removed
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("N2").ClearContents()

# Location info
$ws.Range("P2").Value = "Tveta friluftsgård, 300 m V om, Srm"
$ws.Range("Q2").Value = 648222.682956806
$ws.Range("R2").Value = 6560420.292955686
$ws.Range("S2").Value = 50

# Dates (keep stored as plain text, not auto-converted Excel date serials)
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2019-09-27"
$ws.Range("Y2").Style = "Normal"

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2019-09-27"
$ws.Range("AA2").Style = "Normal"

# Bestamningsmetod (AF2) removed
$ws.Range("AF2").ClearContents()

# Biotop-beskrivning (AI2) new value
$ws.Range("AI2").Value = "barrskog"

# Reporter / observers
$ws.Range("AW2").Value = "Hans Rydberg"
$ws.Range("AX2").Value = "Hans Rydberg"
